# routes.xlsx — "added 2 ltcfs and updated 3 zip code routes"
#
# The "routes" worksheet is a 3-column table (name | zip | region) that maps
# King County zip codes to one of four service regions. This commit:
#   - re-routes the zip on row 50 from South King County  -> East King County
#   - re-routes the zip on row 73 from South King County  -> West King County
#   - re-routes the zip on row 96 from West King County   -> South King County
#   - fixes the region code on row 47 (name stays "South King County",
#     but the machine-readable region key is corrected to east_king_county)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: only the region code changes (name/zip untouched)
$ws.Range("C47").Value = "east_king_county"

# Row 50: re-assigned to East King County (name + region)
$ws.Range("A50").Value = "East King County"
$ws.Range("C50").Value = "east_king_county"

# Row 73: re-assigned to West King County (name + region)
$ws.Range("A73").Value = "West King County"
$ws.Range("C73").Value = "west_king_county"

# Row 96: re-assigned to South King County (name + region)
$ws.Range("A96").Value = "South King County"
$ws.Range("C96").Value = "south_king_county"

# Match the author's final selection/cursor position (cell B1)
$ws.Range("B1").Select()
